# Update random-id string values in the "Custom Table Entry" sheet (first/active sheet)
# per the commit "error fixes and graph report update".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value  = "r287445k397636v"
$ws.Range("A5").Value  = "o0k701899833syn"
$ws.Range("A6").Value  = "29k5z09pu083z17"
$ws.Range("A7").Value  = "ral65f969l2i021"
$ws.Range("A9").Value  = "721775c46nu5dbb"
$ws.Range("A10").Value = "4526v8k0ema390l"
$ws.Range("A11").Value = "6744r3n65e16u49"
$ws.Range("A13").Value = "59iiat6631exp92"
